# Updates cryptocurrency price/volume figures in the worksheet to reflect
# a new data pull (per commit message: "Updated cryptos list ... with
# GitHub Actions"). Only column D (Price) and column E (Volume(1h)) values
# change for the affected rows; coin name, link, row order and styling are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values keyed by A1 reference, taken from the updated data feed.
$updates = [ordered]@{
    'D2' = '59.421.20'
    'E2' = '  -1.08%  '
    'D3' = '2.346.03'
    'E3' = '  -3.03%  '
    'E4' = '  -0.01%  '
    'D5' = '558.74'
    'E5' = '  +1.09%  '
    'D6' = '131.71'
    'E6' = '  -4.04%  '
    'E7' = '  +0.04%  '
    'D8' = '0.577'
    'E8' = '  -2.84%  '
    'E9' = '  -1.63%  '
    'D10' = '5.59'
    'E10' = '  -1.52%  '
    'E11' = '  +1.06%  '
    'E12' = '  -4.17%  '
    'D13' = '24.01'
    'E13' = '  -4.69%  '
    'D14' = '2.768.35'
    'E14' = '  -2.88%  '
    'D15' = '59.403.63'
    'E15' = '  -0.99%  '
    'E16' = '  -1.41%  '
    'D17' = '2.354.25'
    'E17' = '  -2.68%  '
    'D18' = '10.95'
    'E18' = '  -3.07%  '
    'E19' = '  +0.15%  '
    'D20' = '318.98'
    'E20' = '  -2.79%  '
    'E21' = '  -1.05%  '
    'E22' = '  +0.04%  '
    'D23' = '63.82'
    'E23' = '  -3.21%  '
    'E24' = '  -3.10%  '
    'E25' = '  +0.02%  '
    'D26' = '8.34'
    'E26' = '  -3.24%  '
    'E27' = '  -3.18%  '
    'E28' = '  +1.29%  '
    'D29' = '171.08'
    'E29' = '  +1.14%  '
    'D30' = '0.0₃0744'
    'E30' = '  -4.25%  '
    'E31' = '  -1.82%  '
    'E32' = '  +4.72%  '
    'D33' = '0.397'
    'E33' = '  -2.03%  '
    'D34' = '17.98'
    'E34' = '  -3.25%  '
    'E35' = '  +0.02%  '
    'E36' = '  -1.38%  '
    'E37' = '  +0.03%  '
    'D38' = '4.06'
    'E38' = '  -3.08%  '
    'E39' = '  -2.39%  '
    'D40' = '38.51'
    'E40' = '  -2.57%  '
    'D41' = '312.29'
    'E41' = '  -3.74%  '
    'D42' = '144.74'
    'E42' = '  +3.08%  '
    'E43' = '  -5.18%  '
    'E44' = '  -1.65%  '
    'E45' = '  -2.50%  '
    'D46' = '0.564'
    'E46' = '  -2.29%  '
    'D47' = '18.74'
    'E47' = '  -4.54%  '
    'D48' = '0.0215'
    'E48' = '  -3.59%  '
    'E49' = '  +0.23%  '
    'D50' = '4.67'
    'E50' = '  +0.15%  '
    'E51' = '  -0.40%  '
}

# A handful of the new Price values (column D) are plain decimals (e.g.
# "5.59", "0.577"). Assigning those directly would make Excel auto-convert
# the General-formatted cell to a Number, but the sheet stores prices as
# text (other rows, e.g. "59.421.20", use dotted/thousands notation that
# can never be a number). Flag exactly those refs so we can force text
# entry for them while leaving every other assignment as a plain literal.
$forceText = @{
    'D2' = $false
    'E2' = $false
    'D3' = $false
    'E3' = $false
    'E4' = $false
    'D5' = $true
    'E5' = $false
    'D6' = $true
    'E6' = $false
    'E7' = $false
    'D8' = $true
    'E8' = $false
    'E9' = $false
    'D10' = $true
    'E10' = $false
    'E11' = $false
    'E12' = $false
    'D13' = $true
    'E13' = $false
    'D14' = $false
    'E14' = $false
    'D15' = $false
    'E15' = $false
    'E16' = $false
    'D17' = $false
    'E17' = $false
    'D18' = $true
    'E18' = $false
    'E19' = $false
    'D20' = $true
    'E20' = $false
    'E21' = $false
    'E22' = $false
    'D23' = $true
    'E23' = $false
    'E24' = $false
    'E25' = $false
    'D26' = $true
    'E26' = $false
    'E27' = $false
    'E28' = $false
    'D29' = $true
    'E29' = $false
    'D30' = $false
    'E30' = $false
    'E31' = $false
    'E32' = $false
    'D33' = $true
    'E33' = $false
    'D34' = $true
    'E34' = $false
    'E35' = $false
    'E36' = $false
    'E37' = $false
    'D38' = $true
    'E38' = $false
    'E39' = $false
    'D40' = $true
    'E40' = $false
    'D41' = $true
    'E41' = $false
    'D42' = $true
    'E42' = $false
    'E43' = $false
    'E44' = $false
    'E45' = $false
    'D46' = $true
    'E46' = $false
    'D47' = $true
    'E47' = $false
    'D48' = $true
    'E48' = $false
    'E49' = $false
    'D50' = $true
    'E50' = $false
    'E51' = $false
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $cell = $ws.Range($ref)

    if ($forceText[$ref]) {
        # Force text formatting, assign, then restore the default style so
        # no stray number-format/style is left behind on the cell.
        $cell.NumberFormat = '@'
        $cell.Value = $newValue
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $newValue
    }
}
